$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 29) down into the
# new row (row 30) so the new row matches the style of the rest of the table.
$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the values for the new part row
$ws.Range("A30").Value = "Part 029"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "PETG"

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Tabella1")
$table.Resize($ws.Range("A1:C30"))

# Update the view: scroll back to the top and select the cell below the new
# last row, matching the author's final on-screen state.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B31").Select()
